$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    8   = 5
    12  = 2
    16  = 3
    18  = 1
    23  = 1
    24  = 1
    31  = 1
    32  = 2
    33  = 1
    37  = 2
    44  = 5
    45  = 4
    50  = 2
    51  = 1
    54  = 2
    58  = 2
    66  = 2
    72  = 2
    73  = 2
    89  = 2
    91  = 4
    92  = 1
    101 = 2
}

foreach ($row in $changes.Keys) {
    $ws.Range("B$row").Value = $changes[$row]
}
